$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.331.73'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '2.395.80'
$ws.Range('E3').Value = '  -3.74%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'549.00"
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D6').Value = "'142.52"
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'0.541"
$ws.Range('E8').Value = '  -10.42%  '
$ws.Range('D9').Value = '2.395.31'
$ws.Range('E9').Value = '  -3.77%  '
$ws.Range('E10').Value = '  -1.86%  '
$ws.Range('D11').Value = "'0.155"
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = "'5.29"
$ws.Range('E12').Value = '  -2.92%  '
$ws.Range('D13').Value = "'0.349"
$ws.Range('E13').Value = '  -3.23%  '
$ws.Range('D14').Value = "'25.45"
$ws.Range('D15').Value = '2.826.34'
$ws.Range('E15').Value = '  -3.74%  '
$ws.Range('D16').Value = "'0.0000167"
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('D17').Value = '61.195.19'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '2.395.66'
$ws.Range('E18').Value = '  -3.66%  '
$ws.Range('D19').Value = "'10.79"
$ws.Range('E19').Value = '  -3.70%  '
$ws.Range('D20').Value = "'4.16"
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('D21').Value = "'319.58"
$ws.Range('E21').Value = '  -1.26%  '
$ws.Range('D22').Value = "'6.76"
$ws.Range('E22').Value = '  -3.72%  '
$ws.Range('E23').Value = '  +8.44%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = "'63.92"
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('D26').Value = "'8.20"
$ws.Range('E26').Value = '  +8.08%  '
$ws.Range('B27').Value = 'Bittensor'
$ws.Range('C27').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D27').Value = "'542.21"
$ws.Range('E27').Value = '  +0.92%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = "'0.999"
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '0.0₃0946'
$ws.Range('E29').Value = '  -5.69%  '
$ws.Range('D30').Value = '2.511.19'
$ws.Range('E30').Value = '  -3.85%  '
$ws.Range('E31').Value = '  -5.64%  '
$ws.Range('D32').Value = "'8.13"
$ws.Range('E32').Value = '  -3.44%  '
$ws.Range('D33').Value = "'0.147"
$ws.Range('E33').Value = '  -3.42%  '
$ws.Range('D34').Value = "'1.85"
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('D36').Value = "'1.00"
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = "'5.61"
$ws.Range('E37').Value = '  -6.08%  '
$ws.Range('E38').Value = '  -4.12%  '
$ws.Range('E39').Value = '  -2.05%  '
$ws.Range('D40').Value = "'1.86"
$ws.Range('E40').Value = '  +6.15%  '
$ws.Range('D41').Value = "'18.17"
$ws.Range('E41').Value = '  -2.28%  '
$ws.Range('D42').Value = "'139.32"
$ws.Range('E42').Value = '  -6.22%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = "'40.35"
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  -5.34%  '
$ws.Range('D46').Value = "'142.01"
$ws.Range('E46').Value = '  -4.87%  '
$ws.Range('D47').Value = "'3.64"
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').Value = "'20.31"
$ws.Range('E48').Value = '  -4.05%  '
$ws.Range('D49').Value = "'0.0521"
$ws.Range('E49').Value = '  -3.04%  '
$ws.Range('D50').Value = "'0.580"
$ws.Range('E50').Value = '  -3.24%  '
$ws.Range('D51').Value = "'0.0228"
$ws.Range('E51').Value = '  -0.92%  '

# Restore default (unstyled) cell style for text-forced numeric-looking values
# so the saved style matches the original workbook (no explicit style index).
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
